$wb = $excel.ActiveWorkbook

# --- 1) Status text: "Ready for handoff" -> "In Translation" ---------------
# Overview sheet: zh-cn (E2) and de-de (F2) status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-locale sheets: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status columns to match the shorter text ----------------
# (was ~17.22 chars wide for "Ready for handoff"; "In Translation" fits a
# narrower column, ~13.41 chars). ColumnWidth values resolve to Excel's
# internal pixel grid, so 12.5 is the precise input that lands on the
# closest representable width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
